$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Intro paragraph (paragraph #2): replace "Text " with the glossary
#    blurb, give it 12pt (sz/szCs 24 half-points) run + paragraph-mark
#    formatting, and move the "_GoBack" bookmark onto it.
# ---------------------------------------------------------------------

# Remove the old "_GoBack" bookmark first (it currently sits on the
# trailing empty paragraph at the end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$introPara = $d.Paragraphs.Item(2)
$introXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>A list of terms used throughout the document that we defined.</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$introPara.Range.InsertXML($introXml) | Out-Null

# ---------------------------------------------------------------------
# 2. Glossary table edits
# ---------------------------------------------------------------------
$t = $d.Tables.Item(1)

# -- Header row: split "word/acronym" -> "Word/A" + "cronym"
$hdr1 = $t.Rows.Item(1)
$hdr1C1 = $hdr1.Cells.Item(1)
$hdr1C1Xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>Word/A</w:t></w:r><w:r><w:t>cronym</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Range($hdr1C1.Range.Start, $hdr1C1.Range.End - 1).InsertXML($hdr1C1Xml) | Out-Null

# -- Header row: split "definition/meaning" -> "Definition/M" + "eaning"
$hdr1C2 = $hdr1.Cells.Item(2)
$hdr1C2Xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>Definition/M</w:t></w:r><w:r><w:t>eaning</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Range($hdr1C2.Range.Start, $hdr1C2.Range.End - 1).InsertXML($hdr1C2Xml)

# -- Row "FIC(s) / Faculty Information Cards" becomes "CSUS / California
#    State University, Sacramento".
$row3 = $t.Rows.Item(3)
$row3.Cells.Item(1).Range.Text = "CSUS"
$row3.Cells.Item(2).Range.Text = "California State University, Sacramento"

# -- Insert a new "ECS" row ahead of the "Hornet CardGen" row (row 4).
$beforeRow = $t.Rows.Item(4)
$ecsRow = $t.Rows.Add($beforeRow)
$ecsRow.Cells.Item(1).Range.Text = "ECS"
$ecsRow.Cells.Item(2).Range.Text = "Engineering and Computer Science College of CSUS"

# -- Re-insert the original "FIC(s) / Faculty Information Cards" row
#    right after the new ECS row (still ahead of "Hornet CardGen", now
#    row 5).
$beforeRow2 = $t.Rows.Item(5)
$ficRow = $t.Rows.Add($beforeRow2)
$ficRow.Cells.Item(1).Range.Text = "FIC(s)"
$ficRow.Cells.Item(2).Range.Text = "Faculty Information Cards"

# -- Add a new "Registrar" row right after the "Hornet CardGen" row
#    (now row 6), i.e. ahead of the "SRS" row (row 7).
$srsRow = $t.Rows.Item(7)
$registrarRow = $t.Rows.Add($srsRow)
$registrarRow.Cells.Item(1).Range.Text = "Registrar"
$registrarRow.Cells.Item(2).Range.Text = "The CSUS Registrar" + [char]8217 + "s Office which is responsible for courses and other registration concerns."
